$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted cells stay as text (avoid Excel auto-converting
# numeric-looking strings like "1.000" or "0.9999" into numbers).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.715.64'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.80%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.893.92'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.18%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.79'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.62%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9997'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4839'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.50%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2887'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +3.70%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06559'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.849.28'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '16.89'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +5.29%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07462'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.45%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.118'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.90%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '88.19'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6708'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +4.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.700.96'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.96%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.27'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.72%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9997'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.17%  '
$ws.Range('B19').NumberFormat = '@'
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').NumberFormat = '@'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '234.20'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +6.39%  '
$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').NumberFormat = '@'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007590'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.96%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.116.42'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.05%  '
$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.289'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.77%  '
$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = 'BinanceUSD'
$ws.Range('C23').NumberFormat = '@'
$ws.Range('C23').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.201'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +2.46%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '170.27'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +4.19%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.389'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.83'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +3.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.965'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1032'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +13.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.401'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.64%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.355'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +4.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.040'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05082'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.81%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.215'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +7.29%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7532'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +5.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9993'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.714'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.88%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01889'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +3.71%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.646'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9225'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.071'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +2.53%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '107.21'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4306'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.11%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.32%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.652'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.37%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.437'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +3.11%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '64.42'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1281'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.31%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.030'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +4.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '34.09'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.97%  '
